# Update res_bus vm_pu results for the 380 kV case (rows 2-25, cols B-F & I-N)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037195823397943
$ws.Range("D2").Value = 1.058775295629043
$ws.Range("E2").Value = 1.035995151814558
$ws.Range("F2").Value = 1.062890104108977
$ws.Range("I2").Value = 1.046422988364257
$ws.Range("J2").Value = 1.042300013289401
$ws.Range("K2").Value = 1.061506608624412
$ws.Range("L2").Value = 1.038790161008696
$ws.Range("M2").Value = 1.065610226813948
$ws.Range("N2").Value = 1.017949323390593

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038301907098893
$ws.Range("D3").Value = 1.059625727306972
$ws.Range("E3").Value = 1.036941194715247
$ws.Range("F3").Value = 1.063938690463393
$ws.Range("I3").Value = 1.046771927702337
$ws.Range("J3").Value = 1.043049588437844
$ws.Range("K3").Value = 1.062171047564367
$ws.Range("L3").Value = 1.039545351382776
$ws.Range("M3").Value = 1.066473130196653
$ws.Range("N3").Value = 1.018205966331456

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039017611630844
$ws.Range("D4").Value = 1.060175866077687
$ws.Range("E4").Value = 1.037553648142538
$ws.Range("F4").Value = 1.064617505181267
$ws.Range("I4").Value = 1.046996346987633
$ws.Range("J4").Value = 1.043534060973029
$ws.Range("K4").Value = 1.062600191590333
$ws.Range("L4").Value = 1.04003369016254
$ws.Range("M4").Value = 1.067031190109528
$ws.Range("N4").Value = 1.01837165451397

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039318493206447
$ws.Range("D5").Value = 1.060407108717139
$ws.Range("E5").Value = 1.037811195448025
$ws.Range("F5").Value = 1.064902952750207
$ws.Range("I5").Value = 1.047090365359391
$ws.Range("J5").Value = 1.043737601072812
$ws.Range("K5").Value = 1.062780413820245
$ws.Range("L5").Value = 1.040238911271845
$ws.Range("M5").Value = 1.067265727421962
$ws.Range("N5").Value = 1.018441219529261

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039369012532416
$ws.Range("D6").Value = 1.060445933268317
$ws.Range("E6").Value = 1.037854442988265
$ws.Range("F6").Value = 1.064950884979456
$ws.Range("I6").Value = 1.047106132252002
$ws.Range("J6").Value = 1.043771768606456
$ws.Range("K6").Value = 1.062810662782436
$ws.Range("L6").Value = 1.040273364314342
$ws.Range("M6").Value = 1.067305103096639
$ws.Range("N6").Value = 1.018452894504913

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039021632026658
$ws.Range("D7").Value = 1.060178956093951
$ws.Range("E7").Value = 1.037557089219822
$ws.Range("F7").Value = 1.064621319055638
$ws.Range("I7").Value = 1.046997604552604
$ws.Range("J7").Value = 1.04353678120487
$ws.Range("K7").Value = 1.062602600475115
$ws.Range("L7").Value = 1.040036432637605
$ws.Range("M7").Value = 1.067034324288035
$ws.Range("N7").Value = 1.018372584399976

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037569631156617
$ws.Range("D8").Value = 1.059062733266737
$ws.Range("E8").Value = 1.036314809219935
$ws.Range("F8").Value = 1.063244414647948
$ws.Range("I8").Value = 1.046541197234134
$ws.Range("J8").Value = 1.042553450261964
$ws.Range("K8").Value = 1.061731322816439
$ws.Range("L8").Value = 1.039045447408299
$ws.Range("M8").Value = 1.065901910311075
$ws.Range("N8").Value = 1.018036135095609

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035010944492715
$ws.Range("D9").Value = 1.057094694959282
$ws.Range("E9").Value = 1.034128054918739
$ws.Range("F9").Value = 1.060820508369882
$ws.Range("I9").Value = 1.045726475068794
$ws.Range("J9").Value = 1.040816454596851
$ws.Range("K9").Value = 1.060189960565971
$ws.Range("L9").Value = 1.037296751338497
$ws.Range("M9").Value = 1.063904188988144
$ws.Range("N9").Value = 1.017440383467543

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.033305054004772
$ws.Range("D10").Value = 1.055781943906844
$ws.Range("E10").Value = 1.032671764773897
$ws.Range("F10").Value = 1.059206174009831
$ws.Range("I10").Value = 1.045176285201569
$ws.Range("J10").Value = 1.039655592790223
$ws.Range("K10").Value = 1.059158321550698
$ws.Range("L10").Value = 1.036129297039099
$ws.Range("M10").Value = 1.062570858318006
$ws.Range("N10").Value = 1.017041277961602

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.032566349386991
$ws.Range("D11").Value = 1.055213339863933
$ws.Range("E11").Value = 1.032041540509923
$ws.Range("F11").Value = 1.058507529966244
$ws.Range("I11").Value = 1.044936376260656
$ws.Range("J11").Value = 1.03915224360585
$ws.Range("K11").Value = 1.058710647858277
$ws.Range("L11").Value = 1.035623381034277
$ws.Range("M11").Value = 1.061993151926254
$ws.Range("N11").Value = 1.016868001369147

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.032291954067856
$ws.Range("D12").Value = 1.055002109042737
$ws.Range("E12").Value = 1.031807500879393
$ws.Range("F12").Value = 1.058248078466775
$ws.Range("I12").Value = 1.044847011918051
$ws.Range("J12").Value = 1.038965173470439
$ws.Range("K12").Value = 1.058544216618267
$ws.Range("L12").Value = 1.035435400819564
$ws.Range("M12").Value = 1.061778511054641
$ws.Range("N12").Value = 1.016803569394394

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032350813153909
$ws.Range("D13").Value = 1.055047419942517
$ws.Range("E13").Value = 1.031857700746972
$ws.Range("F13").Value = 1.058303729161262
$ws.Range("I13").Value = 1.044866192254752
$ws.Range("J13").Value = 1.039005305355053
$ws.Range("K13").Value = 1.058579923259081
$ws.Range("L13").Value = 1.035475725955617
$ws.Range("M13").Value = 1.061824554740428
$ws.Range("N13").Value = 1.016817393409682

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.032543667944379
$ws.Range("D14").Value = 1.055195879986808
$ws.Range("E14").Value = 1.032022193626364
$ws.Range("F14").Value = 1.058486082482401
$ws.Range("I14").Value = 1.044928994507003
$ws.Range("J14").Value = 1.039136782450754
$ws.Range("K14").Value = 1.05869689356403
$ws.Range("L14").Value = 1.035607843766693
$ws.Range("M14").Value = 1.061975410767296
$ws.Range("N14").Value = 1.016862676816574

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032662491114519
$ws.Range("D15").Value = 1.055287347685705
$ws.Range("E15").Value = 1.032123550230668
$ws.Range("F15").Value = 1.058598443800482
$ws.Range("I15").Value = 1.044967655720353
$ws.Range("J15").Value = 1.039217776056837
$ws.Range("K15").Value = 1.058768943576355
$ws.Range("L15").Value = 1.035689237886752
$ws.Range("M15").Value = 1.062068350854023
$ws.Range("N15").Value = 1.016890568229833

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.033354078401615
$ws.Range("D16").Value = 1.055819676632431
$ws.Range("E16").Value = 1.032713598228004
$ws.Range("F16").Value = 1.059252548573103
$ws.Range("I16").Value = 1.045192171907695
$ws.Range("J16").Value = 1.03968898387916
$ws.Range("K16").Value = 1.059188011802691
$ws.Range("L16").Value = 1.036162864552528
$ws.Range("M16").Value = 1.062609191081575
$ws.Range("N16").Value = 1.017052768037121

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.033787880601
$ws.Range("D17").Value = 1.056153546031135
$ws.Range("E17").Value = 1.033083815708381
$ws.Range("F17").Value = 1.059662951117911
$ws.Range("I17").Value = 1.045332556962701
$ws.Range("J17").Value = 1.03998437550798
$ws.Range("K17").Value = 1.059450623293628
$ws.Range("L17").Value = 1.036459850585795
$ws.Range("M17").Value = 1.062948347968708
$ws.Range("N17").Value = 1.017154388172185

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034040906031086
$ws.Range("D18").Value = 1.056348269583874
$ws.Range("E18").Value = 1.033299791951534
$ws.Range("F18").Value = 1.059902368011538
$ws.Range("I18").Value = 1.045414279768246
$ws.Range("J18").Value = 1.040156606058578
$ws.Range("K18").Value = 1.059603706812652
$ws.Range("L18").Value = 1.036633038883521
$ws.Range("M18").Value = 1.063146137192364
$ws.Range("N18").Value = 1.017213616973514

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034127180510112
$ws.Range("D19").Value = 1.056414662424189
$ws.Range("E19").Value = 1.033373440152831
$ws.Range("F19").Value = 1.059984009099718
$ws.Range("I19").Value = 1.04544211775661
$ws.Range("J19").Value = 1.040215320977513
$ws.Range("K19").Value = 1.059655888507441
$ws.Range("L19").Value = 1.036692085074135
$ws.Range("M19").Value = 1.063213572230105
$ws.Range("N19").Value = 1.017233804916445

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033741338157538
$ws.Range("D20").Value = 1.056117726742431
$ws.Range("E20").Value = 1.033044091293779
$ws.Range("F20").Value = 1.059618915082361
$ws.Range("I20").Value = 1.045317511672281
$ws.Range("J20").Value = 1.039952689651573
$ws.Range("K20").Value = 1.059422457218384
$ws.Range("L20").Value = 1.0364279907929
$ws.Range("M20").Value = 1.062911963291167
$ws.Range("N20").Value = 1.017143489896883

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032486877232194
$ws.Range("D21").Value = 1.055152162919092
$ws.Range("E21").Value = 1.031973753091969
$ws.Range("F21").Value = 1.05843238242908
$ws.Range("I21").Value = 1.044910507752028
$ws.Range("J21").Value = 1.039098068592809
$ws.Range("K21").Value = 1.058662452735479
$ws.Range("L21").Value = 1.035568940036366
$ws.Range("M21").Value = 1.061930988946289
$ws.Range("N21").Value = 1.016849343896094

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031698103727694
$ws.Range("D22").Value = 1.054544924508758
$ws.Range("E22").Value = 1.031301099696341
$ws.Range("F22").Value = 1.057686686240796
$ws.Range("I22").Value = 1.044653153190723
$ws.Range("J22").Value = 1.038560133439428
$ws.Range("K22").Value = 1.058183767108315
$ws.Range("L22").Value = 1.035028470636142
$ws.Range("M22").Value = 1.061313893206674
$ws.Range("N22").Value = 1.016664001474727

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.032116252038642
$ws.Range("D23").Value = 1.054866847240404
$ws.Range("E23").Value = 1.031657656606879
$ws.Range("F23").Value = 1.058081963220248
$ws.Range("I23").Value = 1.044789719672834
$ws.Range("J23").Value = 1.038845360147313
$ws.Range("K23").Value = 1.058437607096149
$ws.Range("L23").Value = 1.035315016976489
$ws.Range("M23").Value = 1.061641057531408
$ws.Range("N23").Value = 1.01676229303468

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033762368696515
$ws.Range("D24").Value = 1.056133911989048
$ws.Range("E24").Value = 1.033062040939589
$ws.Range("F24").Value = 1.059638812959822
$ws.Range("I24").Value = 1.045324310490069
$ws.Range("J24").Value = 1.039967007331114
$ws.Range("K24").Value = 1.059435184543622
$ws.Range("L24").Value = 1.036442386981111
$ws.Range("M24").Value = 1.062928404069555
$ws.Range("N24").Value = 1.017148414495842

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035672440291702
$ws.Range("D25").Value = 1.057603609475648
$ws.Range("E25").Value = 1.034693110702603
$ws.Range("F25").Value = 1.061446864123486
$ws.Range("I25").Value = 1.045938341469616
$ws.Range("J25").Value = 1.041266013554785
$ws.Range("K25").Value = 1.060589156720328
$ws.Range("L25").Value = 1.037749122692682
$ws.Range("M25").Value = 1.064420915495446
$ws.Range("N25").Value = 1.017594741227114
